$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-26 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-27 Saturday", 2) | Out-Null
$d.Content.Find.Execute("96-22=", $true, $false, $false, $false, $false, $true, 1, $false, "33+10=", 2) | Out-Null
$d.Content.Find.Execute("94-11=", $true, $false, $false, $false, $false, $true, 1, $false, "10+64=", 2) | Out-Null
$d.Content.Find.Execute("16+40=", $true, $false, $false, $false, $false, $true, 1, $false, "37+27=", 2) | Out-Null
$d.Content.Find.Execute("72+0=", $true, $false, $false, $false, $false, $true, 1, $false, "16+9=", 2) | Out-Null
$d.Content.Find.Execute("62-50=", $true, $false, $false, $false, $false, $true, 1, $false, "61-15=", 2) | Out-Null
$d.Content.Find.Execute("93-46=", $true, $false, $false, $false, $false, $true, 1, $false, "39-14=", 2) | Out-Null
$d.Content.Find.Execute("95-57=", $true, $false, $false, $false, $false, $true, 1, $false, "39+38=", 2) | Out-Null
$d.Content.Find.Execute("24+3=", $true, $false, $false, $false, $false, $true, 1, $false, "56+1=", 2) | Out-Null
$d.Content.Find.Execute("90-49=", $true, $false, $false, $false, $false, $true, 1, $false, "56+19=", 2) | Out-Null
$d.Content.Find.Execute("50-11=", $true, $false, $false, $false, $false, $true, 1, $false, "1+31=", 2) | Out-Null
$d.Content.Find.Execute("93-77=", $true, $false, $false, $false, $false, $true, 1, $false, "88-41=", 2) | Out-Null
$d.Content.Find.Execute("38+57=", $true, $false, $false, $false, $false, $true, 1, $false, "76-12=", 2) | Out-Null
$d.Content.Find.Execute("43-19=", $true, $false, $false, $false, $false, $true, 1, $false, "85-82=", 2) | Out-Null
$d.Content.Find.Execute("47-26=", $true, $false, $false, $false, $false, $true, 1, $false, "93-58=", 2) | Out-Null
$d.Content.Find.Execute("93-49=", $true, $false, $false, $false, $false, $true, 1, $false, "7+81=", 2) | Out-Null
$d.Content.Find.Execute("15-11=", $true, $false, $false, $false, $false, $true, 1, $false, "30+11=", 2) | Out-Null
$d.Content.Find.Execute("17+64=", $true, $false, $false, $false, $false, $true, 1, $false, "37-30=", 2) | Out-Null
$d.Content.Find.Execute("21+19=", $true, $false, $false, $false, $false, $true, 1, $false, "87-33=", 2) | Out-Null
$d.Content.Find.Execute("39-28=", $true, $false, $false, $false, $false, $true, 1, $false, "1+79=", 2) | Out-Null
$d.Content.Find.Execute("28+69=", $true, $false, $false, $false, $false, $true, 1, $false, "12+75=", 2) | Out-Null
$d.Content.Find.Execute("0+46=", $true, $false, $false, $false, $false, $true, 1, $false, "96-92=", 2) | Out-Null
$d.Content.Find.Execute("66-32=", $true, $false, $false, $false, $false, $true, 1, $false, "50-16=", 2) | Out-Null
$d.Content.Find.Execute("47-21=", $true, $false, $false, $false, $false, $true, 1, $false, "98-90=", 2) | Out-Null
$d.Content.Find.Execute("56+3=", $true, $false, $false, $false, $false, $true, 1, $false, "81-71=", 2) | Out-Null
$d.Content.Find.Execute("99-14=", $true, $false, $false, $false, $false, $true, 1, $false, "2+91=", 2) | Out-Null
$d.Content.Find.Execute("97-73=", $true, $false, $false, $false, $false, $true, 1, $false, "47-15=", 2) | Out-Null
$d.Content.Find.Execute("45-39=", $true, $false, $false, $false, $false, $true, 1, $false, "51-22=", 2) | Out-Null
$d.Content.Find.Execute("19+44=", $true, $false, $false, $false, $false, $true, 1, $false, "46+46=", 2) | Out-Null
$d.Content.Find.Execute("34-15=", $true, $false, $false, $false, $false, $true, 1, $false, "67-5=", 2) | Out-Null
$d.Content.Find.Execute("79-27=", $true, $false, $false, $false, $false, $true, 1, $false, "95-31=", 2) | Out-Null
$d.Content.Find.Execute("74-25=", $true, $false, $false, $false, $false, $true, 1, $false, "89-62=", 2) | Out-Null
$d.Content.Find.Execute("29+60=", $true, $false, $false, $false, $false, $true, 1, $false, "93-60=", 2) | Out-Null
$d.Content.Find.Execute("86-44=", $true, $false, $false, $false, $false, $true, 1, $false, "69+30=", 2) | Out-Null
$d.Content.Find.Execute("54-49=", $true, $false, $false, $false, $false, $true, 1, $false, "57-36=", 2) | Out-Null
$d.Content.Find.Execute("15+72=", $true, $false, $false, $false, $false, $true, 1, $false, "52-17=", 2) | Out-Null
$d.Content.Find.Execute("32-23=", $true, $false, $false, $false, $false, $true, 1, $false, "45-16=", 2) | Out-Null
$d.Content.Find.Execute("5+25=", $true, $false, $false, $false, $false, $true, 1, $false, "11+16=", 2) | Out-Null
$d.Content.Find.Execute("12-0=", $true, $false, $false, $false, $false, $true, 1, $false, "59+17=", 2) | Out-Null
$d.Content.Find.Execute("4+49=", $true, $false, $false, $false, $false, $true, 1, $false, "3+69=", 2) | Out-Null
$d.Content.Find.Execute("56-30=", $true, $false, $false, $false, $false, $true, 1, $false, "62-60=", 2) | Out-Null
$d.Content.Find.Execute("18-4=", $true, $false, $false, $false, $false, $true, 1, $false, "33-5=", 2) | Out-Null
$d.Content.Find.Execute("87-17=", $true, $false, $false, $false, $false, $true, 1, $false, "39-3=", 2) | Out-Null
$d.Content.Find.Execute("78-6=", $true, $false, $false, $false, $false, $true, 1, $false, "20+53=", 2) | Out-Null
$d.Content.Find.Execute("3+64=", $true, $false, $false, $false, $false, $true, 1, $false, "42+43=", 2) | Out-Null
$d.Content.Find.Execute("74-43=", $true, $false, $false, $false, $false, $true, 1, $false, "72-60=", 2) | Out-Null
$d.Content.Find.Execute("4+0=", $true, $false, $false, $false, $false, $true, 1, $false, "92-33=", 2) | Out-Null
$d.Content.Find.Execute("91-54=", $true, $false, $false, $false, $false, $true, 1, $false, "10+25=", 2) | Out-Null
$d.Content.Find.Execute("41+22=", $true, $false, $false, $false, $false, $true, 1, $false, "77+1=", 2) | Out-Null
$d.Content.Find.Execute("44+31=", $true, $false, $false, $false, $false, $true, 1, $false, "6+78=", 2) | Out-Null
$d.Content.Find.Execute("99-36=", $true, $false, $false, $false, $false, $true, 1, $false, "72-67=", 2) | Out-Null
$d.Content.Find.Execute("54-48=", $true, $false, $false, $false, $false, $true, 1, $false, "94-25=", 2) | Out-Null
$d.Content.Find.Execute("17+78=", $true, $false, $false, $false, $false, $true, 1, $false, "61-32=", 2) | Out-Null
$d.Content.Find.Execute("8+19=", $true, $false, $false, $false, $false, $true, 1, $false, "69-27=", 2) | Out-Null
$d.Content.Find.Execute("54+35=", $true, $false, $false, $false, $false, $true, 1, $false, "57+24=", 2) | Out-Null
$d.Content.Find.Execute("72-29=", $true, $false, $false, $false, $false, $true, 1, $false, "43+36=", 2) | Out-Null
$d.Content.Find.Execute("85-28=", $true, $false, $false, $false, $false, $true, 1, $false, "76+21=", 2) | Out-Null
$d.Content.Find.Execute("39+44=", $true, $false, $false, $false, $false, $true, 1, $false, "81-8=", 2) | Out-Null
$d.Content.Find.Execute("52-1=", $true, $false, $false, $false, $false, $true, 1, $false, "98-13=", 2) | Out-Null
$d.Content.Find.Execute("24+63=", $true, $false, $false, $false, $false, $true, 1, $false, "29-27=", 2) | Out-Null
$d.Content.Find.Execute("43-8=", $true, $false, $false, $false, $false, $true, 1, $false, "60-36=", 2) | Out-Null
$d.Content.Find.Execute("30+3=", $true, $false, $false, $false, $false, $true, 1, $false, "25-16=", 2) | Out-Null
$d.Content.Find.Execute("61+1=", $true, $false, $false, $false, $false, $true, 1, $false, "75+21=", 2) | Out-Null
$d.Content.Find.Execute("80-36=", $true, $false, $false, $false, $false, $true, 1, $false, "2+93=", 2) | Out-Null
$d.Content.Find.Execute("84+13=", $true, $false, $false, $false, $false, $true, 1, $false, "41+31=", 2) | Out-Null
$d.Content.Find.Execute("32-13=", $true, $false, $false, $false, $false, $true, 1, $false, "86-83=", 2) | Out-Null
$d.Content.Find.Execute("74-2=", $true, $false, $false, $false, $false, $true, 1, $false, "1+94=", 2) | Out-Null
$d.Content.Find.Execute("3+49=", $true, $false, $false, $false, $false, $true, 1, $false, "55-24=", 2) | Out-Null
$d.Content.Find.Execute("50+20=", $true, $false, $false, $false, $false, $true, 1, $false, "75-41=", 2) | Out-Null
$d.Content.Find.Execute("14-12=", $true, $false, $false, $false, $false, $true, 1, $false, "77-9=", 2) | Out-Null
$d.Content.Find.Execute("7+38=", $true, $false, $false, $false, $false, $true, 1, $false, "11+77=", 2) | Out-Null
$d.Content.Find.Execute("0+31=", $true, $false, $false, $false, $false, $true, 1, $false, "10+24=", 2) | Out-Null
$d.Content.Find.Execute("30+50=", $true, $false, $false, $false, $false, $true, 1, $false, "62-38=", 2) | Out-Null
$d.Content.Find.Execute("27-15=", $true, $false, $false, $false, $false, $true, 1, $false, "5+16=", 2) | Out-Null
$d.Content.Find.Execute("23+11=", $true, $false, $false, $false, $false, $true, 1, $false, "13-3=", 2) | Out-Null
$d.Content.Find.Execute("29+11=", $true, $false, $false, $false, $false, $true, 1, $false, "77+20=", 2) | Out-Null
$d.Content.Find.Execute("10+17=", $true, $false, $false, $false, $false, $true, 1, $false, "28+41=", 2) | Out-Null
$d.Content.Find.Execute("43+37=", $true, $false, $false, $false, $false, $true, 1, $false, "1+61=", 2) | Out-Null
$d.Content.Find.Execute("84-19=", $true, $false, $false, $false, $false, $true, 1, $false, "29+19=", 2) | Out-Null
$d.Content.Find.Execute("54+28=", $true, $false, $false, $false, $false, $true, 1, $false, "45-27=", 2) | Out-Null
$d.Content.Find.Execute("93-28=", $true, $false, $false, $false, $false, $true, 1, $false, "81+9=", 2) | Out-Null
$d.Content.Find.Execute("18+62=", $true, $false, $false, $false, $false, $true, 1, $false, "51+45=", 2) | Out-Null
$d.Content.Find.Execute("27-1=", $true, $false, $false, $false, $false, $true, 1, $false, "30-1=", 2) | Out-Null
$d.Content.Find.Execute("94-27=", $true, $false, $false, $false, $false, $true, 1, $false, "37-27=", 2) | Out-Null
$d.Content.Find.Execute("85-38=", $true, $false, $false, $false, $false, $true, 1, $false, "38-37=", 2) | Out-Null
$d.Content.Find.Execute("98-46=", $true, $false, $false, $false, $false, $true, 1, $false, "90-21=", 2) | Out-Null
$d.Content.Find.Execute("11+37=", $true, $false, $false, $false, $false, $true, 1, $false, "82+1=", 2) | Out-Null
$d.Content.Find.Execute("3+57=", $true, $false, $false, $false, $false, $true, 1, $false, "94-54=", 2) | Out-Null
$d.Content.Find.Execute("8+6=", $true, $false, $false, $false, $false, $true, 1, $false, "12+50=", 2) | Out-Null
$d.Content.Find.Execute("20+60=", $true, $false, $false, $false, $false, $true, 1, $false, "51+24=", 2) | Out-Null
$d.Content.Find.Execute("97-86=", $true, $false, $false, $false, $false, $true, 1, $false, "74-50=", 2) | Out-Null
$d.Content.Find.Execute("62-10=", $true, $false, $false, $false, $false, $true, 1, $false, "39+22=", 2) | Out-Null
$d.Content.Find.Execute("73-26=", $true, $false, $false, $false, $false, $true, 1, $false, "49+33=", 2) | Out-Null
$d.Content.Find.Execute("0+94=", $true, $false, $false, $false, $false, $true, 1, $false, "65-61=", 2) | Out-Null
$d.Content.Find.Execute("16+79=", $true, $false, $false, $false, $false, $true, 1, $false, "31+47=", 2) | Out-Null
$d.Content.Find.Execute("80+6=", $true, $false, $false, $false, $false, $true, 1, $false, "38-8=", 2) | Out-Null
$d.Content.Find.Execute("10+61=", $true, $false, $false, $false, $false, $true, 1, $false, "21+59=", 2) | Out-Null
$d.Content.Find.Execute("5+79=", $true, $false, $false, $false, $false, $true, 1, $false, "63-16=", 2) | Out-Null
$d.Content.Find.Execute("25+40=", $true, $false, $false, $false, $false, $true, 1, $false, "27+11=", 2) | Out-Null
$d.Content.Find.Execute("65+18=", $true, $false, $false, $false, $false, $true, 1, $false, "86-11=", 2) | Out-Null
$d.Content.Find.Execute("68+15=", $true, $false, $false, $false, $false, $true, 1, $false, "85-29=", 2) | Out-Null
